# producto 6 version preliminar
# Adds a new worksheet "tabla03j" (after the existing "tabla02j" sheet) with a
# small table: anio / ingresado / resuelto / prop_resuelto / ejecutado / prop_ejecutado

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after the current last sheet (tabla02j) so it
# lands at the end of the tab strip, matching the source workbook order.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "tabla03j"

# --- Header row -----------------------------------------------------------
$headers = @("anio", "ingresado", "resuelto", "prop_resuelto", "ejecutado", "prop_ejecutado")
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}
$headerRange = $ws.Range("A1:F1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter

# --- Data rows --------------------------------------------------------------
# anio, ingresado, resuelto, prop_resuelto, ejecutado, prop_ejecutado
$data = @(
    @("2019", 34, 33, 97.1, 18, 52.9),
    @("2020", 19, 19, 100,  10, 52.6),
    @("2021", 31, 29, 93.5, 12, 38.7),
    @("2022", 44, 41, 93.2, 31, 70.5),
    @("2023", 67, 56, 83.6, 31, 46.3),
    @("2024", 49, 18, 36.7, 6,  12.2)
)

$row = 2
foreach ($entry in $data) {
    # "anio" is stored as text (matches the other tabla0*j sheets), so force
    # text typing via NumberFormat, then strip the format back off so no
    # residual numFmt/style survives on the cell.
    $yearCell = $ws.Cells.Item($row, 1)
    $yearCell.NumberFormat = "@"
    $yearCell.Value = $entry[0]
    $yearCell.ClearFormats()

    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 6).Value = $entry[5]

    $row++
}
